{"js": "// \"Inclus\u00e3o do execut\u00e1vel do quali.\"\n// 1) Drop the \" (Suspeito Detido)\" suffix from the title row.\n// 2) Shrink the \"{Imagem}\" placeholder (paragraph mark + run) to 6pt (sz/szCs=12).\n// 3) Collapse the four single-character runs (\"N\",\"7\",\"D\",\"8\") of the secondary\n//    address cell into one run containing \"N7D8\".\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// --- 1. Title: \"DADOS GERAIS DE QUALIFICA\u00c7\u00c3O (Suspeito Detido)\" -> \"DADOS GERAIS DE QUALIFICA\u00c7\u00c3O\"\nconst titleCell = table.getCell(0, 0);\nconst suffix = titleCell.body.search(\" (Suspeito Detido)\", { matchCase: true });\nsuffix.load(\"items\");\nawait context.sync();\nif (suffix.items.length > 0) {\n  suffix.items[0].delete();\n  await context.sync();\n}\n\n// --- 2. \"{Imagem}\" placeholder paragraph: add sz/szCs=12 (6pt) to the paragraph\n//        mark properties AND to the run itself, preserving everything else.\nconst imageCell = table.getCell(1, 1);\nimageCell.body.paragraphs.load(\"items\");\nawait context.sync();\nconst imagePara = imageCell.body.paragraphs.items[0];\n\nconst wordNs =\n  'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' +\n  'xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"';\nconst imageOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document ' + wordNs + '>' +\n  '<w:body>' +\n  '<w:p w14:paraId=\"0B7D2644\" w14:textId=\"6479FCD2\" w:rsidR=\"001A72F2\" w:rsidRPr=\"00D63B0A\" w:rsidRDefault=\"00F86F71\" w:rsidP=\"00561B0C\">' +\n  '<w:pPr><w:jc w:val=\"center\"/><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"12\"/><w:szCs w:val=\"12\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:noProof/><w:sz w:val=\"12\"/><w:szCs w:val=\"12\"/></w:rPr><w:t>{Imagem}</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nimagePara.insertOoxml(imageOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 3. \"Endere\u00e7o(s) secund\u00e1rio(s)\" value cell: merge N + 7 + D + 8 runs into \"N7D8\"\nconst addressCell = table.getCell(10, 1);\naddressCell.body.paragraphs.load(\"items\");\nawait context.sync();\nconst addressPara = addressCell.body.paragraphs.items[0];\nconst addressRange = addressPara.getRange();\naddressRange.insertText(\"N7D8\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# \"Inclus\u00e3o do execut\u00e1vel do quali.\"\n# 1) Drop the \" (Suspeito Detido)\" suffix from the title row.\n# 2) Shrink the \"{Imagem}\" placeholder (paragraph mark + run) to 6pt (sz/szCs=12).\n# 3) Collapse the four single-character runs (\"N\",\"7\",\"D\",\"8\") of the secondary\n#    address cell into one run containing \"N7D8\".\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# wdReplaceAll = 2\n$wdReplaceAll = 2\n\n# --- 1. Title: \"DADOS GERAIS DE QUALIFICA\u00c7\u00c3O (Suspeito Detido)\" -> \"DADOS GERAIS DE QUALIFICA\u00c7\u00c3O\"\n$titleRange = $t.Cell(1, 1).Range\n$titleRange.Find.Execute(\" (Suspeito Detido)\", $true, $false, $false, $false, $false, $true, 1, $false, \"\", $wdReplaceAll) | Out-Null\n\n# --- 2. \"{Imagem}\" placeholder: set size to 6pt (sz=12) and complex-script size\n#        to 6pt (szCs=12). The cell's Range spans the run AND the paragraph\n#        mark, so this updates both the run's rPr and the paragraph's pPr/rPr.\n$imageRange = $t.Cell(2, 2).Range\n$imageRange.Font.Size = 6\n$imageRange.Font.SizeBi = 6\n\n# --- 3. \"Endere\u00e7o(s) secund\u00e1rio(s)\" value cell: merge N + 7 + D + 8 runs into \"N7D8\"\n$addressRange = $t.Cell(11, 3).Range\n$addressRange.Find.Execute(\"N7D8\", $true, $false, $false, $false, $false, $true, 1, $false, \"N7D8\", $wdReplaceAll) | Out-Null\n"}
